$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above existing row 2 ("Rule3") to make room for Rule1 and Rule2,
# pushing the existing Rule3 row down to row 4.
$ws.Rows("2:3").Insert()

# Row 2: Rule1
$ws.Range("A2").Value = "Rule1"
$ws.Range("B2").Value = "GL, CC, PC, Reporting Id"
$ws.Range("C2").Value = "GL, CC, PC, Reporting Id"
$ws.Range("D2").Value = "Exact"

# Row 3: Rule2
$ws.Range("A3").Value = "Rule2"
$ws.Range("B3").Value = "ISIN"
$ws.Range("C3").Value = "ISIN"
$ws.Range("D3").Value = "Exact"

# Row 4: Rule3 (already has RuleName from before; fill in remaining attributes)
$ws.Range("B4").Value = "ISIN, Period"
$ws.Range("C4").Value = "ISIN, Period"
$ws.Range("D4").Value = "Exact"
